$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "60.306.18"
$ws.Range("E2").Value = "  -4.26%  "
$ws.Range("D3").Value = "2.986.12"
$ws.Range("E3").Value = "  -5.84%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.63"
$ws.Range("E5").Value = "  -2.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "125.24"
$ws.Range("E6").Value = "  -6.70%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "2.982.65"
$ws.Range("E8").Value = "  -5.94%  "
$ws.Range("E9").Value = "  -2.96%  "
$ws.Range("E10").Value = "  -5.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.10"
$ws.Range("E11").Value = "  -1.88%  "
$ws.Range("E12").Value = "  -2.83%  "
$ws.Range("E13").Value = "  -5.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.47"
$ws.Range("E14").Value = "  -5.11%  "
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").Value = "3.475.39"
$ws.Range("E16").Value = "  -5.91%  "
$ws.Range("D17").Value = "2.984.17"
$ws.Range("E17").Value = "  -5.80%  "
$ws.Range("D18").Value = "60.200.06"
$ws.Range("E18").Value = "  -4.44%  "
$ws.Range("E19").Value = "  -4.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "431.63"
$ws.Range("E20").Value = "  -5.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.06"
$ws.Range("E21").Value = "  -6.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.660"
$ws.Range("E22").Value = "  -4.73%  "
$ws.Range("E23").Value = "  -7.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.69"
$ws.Range("E24").Value = "  -3.96%  "
$ws.Range("E25").Value = "  -3.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  -4.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.28"
$ws.Range("E29").Value = "  -4.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.88"
$ws.Range("E30").Value = "  -6.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.13"
$ws.Range("E31").Value = "  -8.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.32"
$ws.Range("E32").Value = "  -6.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0931"
$ws.Range("E33").Value = "  -8.06%  "
$ws.Range("E34").Value = "  -7.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.951"
$ws.Range("E35").Value = "  -7.51%  "
$ws.Range("E36").Value = "  -3.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "49.52"
$ws.Range("E37").Value = "  -3.01%  "
$ws.Range("D38").Value = "0.0₃0658"
$ws.Range("E38").Value = "  -6.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.01"
$ws.Range("E39").Value = "  -0.76%  "
$ws.Range("E40").Value = "  -6.75%  "
$ws.Range("E41").Value = "  -1.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "383.96"
$ws.Range("E43").Value = "  -7.03%  "
$ws.Range("D44").Value = "2.630.81"
$ws.Range("E44").Value = "  -6.32%  "
$ws.Range("E46").Value = "  -5.88%  "
$ws.Range("E47").Value = "  -4.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "118.36"
$ws.Range("E48").Value = "  -4.88%  "
$ws.Range("E49").Value = "  -3.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.53"
$ws.Range("E50").Value = "  -5.99%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.01"
$ws.Range("E51").Value = "  -5.33%  "
